$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'62.853.05"
$ws.Range("E2").Value = "'  +1.53%  "

# Row 3
$ws.Range("D3").Value = "'2.445.51"
$ws.Range("E3").Value = "'  +1.85%  "

# Row 4
$ws.Range("E4").Value = "'  +0.14%  "

# Row 5
$ws.Range("D5").Value = "'567.37"
$ws.Range("E5").Value = "'  +0.97%  "

# Row 6
$ws.Range("D6").Value = "'146.13"
$ws.Range("E6").Value = "'  +2.75%  "

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "'  -0.11%  "

# Row 8
$ws.Range("D8").Value = "'0.534"
$ws.Range("E8").Value = "'  +0.47%  "

# Row 9
$ws.Range("D9").Value = "'0.112"
$ws.Range("E9").Value = "'  +2.87%  "

# Row 10
$ws.Range("E10").Value = "'  +0.58%  "

# Row 11
$ws.Range("D11").Value = "'5.32"
$ws.Range("E11").Value = "'  +1.32%  "

# Row 12
$ws.Range("D12").Value = "'0.355"
$ws.Range("E12").Value = "'  +1.78%  "

# Row 13
$ws.Range("D13").Value = "'27.05"
$ws.Range("E13").Value = "'  +5.95%  "

# Row 14
$ws.Range("D14").Value = "'0.0000183"
$ws.Range("E14").Value = "'  +6.50%  "

# Row 15
$ws.Range("D15").Value = "'2.799.42"
$ws.Range("E15").Value = "'  -1.25%  "

# Row 16
$ws.Range("D16").Value = "'62.615.76"
$ws.Range("E16").Value = "'  +1.22%  "

# Row 17
$ws.Range("D17").Value = "'2.448.05"
$ws.Range("E17").Value = "'  +1.70%  "

# Row 18
$ws.Range("D18").Value = "'11.30"
$ws.Range("E18").Value = "'  +0.80%  "

# Row 19
$ws.Range("D19").Value = "'6.96"
$ws.Range("E19").Value = "'  +2.56%  "

# Row 20
$ws.Range("D20").Value = "'324.35"
$ws.Range("E20").Value = "'  +0.99%  "

# Row 21
$ws.Range("D21").Value = "'4.18"
$ws.Range("E21").Value = "'  +1.19%  "

# Row 22
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "'  -0.06%  "

# Row 23
$ws.Range("E23").Value = "'  +6.17%  "

# Row 24
$ws.Range("D24").Value = "'67.45"
$ws.Range("E24").Value = "'  +2.21%  "

# Row 25
$ws.Range("D25").Value = "'8.72"
$ws.Range("E25").Value = "'  -0.67%  "

# Row 26
$ws.Range("B26").Value = "Bittensor"
$ws.Range("C26").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D26").Value = "'588.15"
$ws.Range("E26").Value = "'  +5.05%  "

# Row 27
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "'0.0000102"
$ws.Range("E27").Value = "'  +9.55%  "

# Row 28
$ws.Range("D28").Value = "'2.563.78"
$ws.Range("E28").Value = "'  +1.73%  "

# Row 29
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'8.48"
$ws.Range("E29").Value = "'  +4.10%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "'  -0.24%  "

# Row 31
$ws.Range("D31").Value = "'1.46"
$ws.Range("E31").Value = "'  +4.88%  "

# Row 32
$ws.Range("D32").Value = "'0.148"
$ws.Range("E32").Value = "'  +0.81%  "

# Row 33
$ws.Range("E33").Value = "'  +0.55%  "

# Row 34
$ws.Range("E34").Value = "'  +3.70%  "

# Row 35
$ws.Range("D35").Value = "'4.87"
$ws.Range("E35").Value = "'  +4.78%  "

# Row 36
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "'  -0.14%  "

# Row 37
$ws.Range("E37").Value = "'  +1.62%  "

# Row 38
$ws.Range("D38").Value = "'5.47"
$ws.Range("E38").Value = "'  +1.04%  "

# Row 39
$ws.Range("D39").Value = "'18.84"
$ws.Range("E39").Value = "'  +1.57%  "

# Row 40
$ws.Range("D40").Value = "'148.04"
$ws.Range("E40").Value = "'  -2.62%  "

# Row 41
$ws.Range("D41").Value = "'1.83"
$ws.Range("E41").Value = "'  +2.34%  "

# Row 42
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "'  +0.50%  "

# Row 43
$ws.Range("D43").Value = "'2.44"
$ws.Range("E43").Value = "'  +8.63%  "

# Row 44
$ws.Range("D44").Value = "'149.63"
$ws.Range("E44").Value = "'  +1.73%  "

# Row 45
$ws.Range("D45").Value = "'3.68"
$ws.Range("E45").Value = "'  +2.40%  "

# Row 46
$ws.Range("D46").Value = "'0.0537"
$ws.Range("E46").Value = "'  +1.53%  "

# Row 47
$ws.Range("D47").Value = "'20.61"
$ws.Range("E47").Value = "'  +4.42%  "

# Row 48
$ws.Range("E48").Value = "'  +2.97%  "

# Row 49
$ws.Range("D49").Value = "'0.0232"
$ws.Range("E49").Value = "'  +3.64%  "

# Row 50
$ws.Range("D50").Value = "'0.0926"
$ws.Range("E50").Value = "'  +1.03%  "

# Row 51
$ws.Range("E51").Value = "'  +3.48%  "
